$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1, columns E:BL currently hold text labels like "1960 [YR1960]" ... "2019 [YR2019]".
# Replace them with plain numeric year values and left-align them.
for ($i = 0; $i -lt 60; $i++) {
    $col = 5 + $i
    $year = 1960 + $i
    $ws.Cells.Item(1, $col).Value = $year
}
$ws.Range("E1:BL1").HorizontalAlignment = -4131

# Update the active selection / scroll position to match the edited range.
$ws.Range("E1:BL1").Select() | Out-Null
